$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8046784400939941
$ws.Range("B1").Value = 1.829840540885925
$ws.Range("C1").Value = 4.528347492218018
$ws.Range("D1").Value = 1.430051684379578
$ws.Range("E1").Value = 1.451180100440979
